$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '42.976.27'
$ws.Range("E2").Value = '  +0.50%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.302.79'
$ws.Range("E3").Value = '  +0.18%  '

# Row 4: TetherUSD
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '

# Row 5: BNB
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '302.00'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.33%  '

# Row 6: Solana
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '100.65'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.72%  '

# Row 7: XRP
$ws.Range("E7").Value = '  +0.07%  '

# Row 8: USDC
$ws.Range("E8").Value = '  +0.06%  '

# Row 9: Cardano
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.514'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +3.64%  '

# Row 10: Avalanche
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '36.12'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +7.13%  '

# Row 11: Dogecoin
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0791'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.07%  '

# Row 12: Chainlink
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '18.58'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +12.67%  '

# Row 13: TRON
$ws.Range("E13").Value = '  +1.96%  '

# Row 14: Polkadot
$ws.Range("E14").Value = '  +2.74%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '2.661.85'
$ws.Range("E15").Value = '  +0.23%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '2.326.74'
$ws.Range("E16").Value = '  +0.61%  '

# Row 17: Polygon
$ws.Range("E17").Value = '  +0.55%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '42.867.08'
$ws.Range("E18").Value = '  +0.51%  '

# Row 19: InternetComputer(DFINITY)
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '12.40'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +6.17%  '

# Row 20: Uniswap
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.23'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +3.57%  '

# Row 21: ShibaInu
$ws.Range("D21").Value = '0.0₃0901'
$ws.Range("E21").Value = '  +0.33%  '

# Row 22: Litecoin
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '68.05'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +1.72%  '

# Row 23: ImmutableX
$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '2.27'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +14.35%  '

# Row 24: BitcoinCash
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '236.40'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.07%  '

# Row 25: Dai
$ws.Range("E25").Value = '  +0.30%  '

# Row 26: PancakeSwap
$ws.Range("E26").Value = '  -0.37%  '

# Row 27: EthereumClassic
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '24.81'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +3.11%  '

# Row 28: Toncoin
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.29'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +2.15%  '

# Row 29: Monero
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '169.83'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.45%  '

# Row 30: InjectiveProtocol
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '34.64'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.57%  '

# Row 31: Cosmos
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '9.18'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.16%  '

# Row 32: FirstDigitalUSD
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.07%  '

# Row 33: Filecoin
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.02'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.88%  '

# Row 34: Celestia
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '17.80'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +3.74%  '

# Row 35: RenderToken
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '4.70'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.67%  '

# Row 36: WEMIXToken
$ws.Range("E36").Value = '  +1.94%  '

# Row 37: Hedera
$ws.Range("E37").Value = '  -0.40%  '

# Row 38: LidoDAOToken
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.84'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.69%  '

# Row 39: Kaspa
$ws.Range("E39").Value = '  +1.66%  '

# Row 40: ARBITRUM
$ws.Range("E40").Value = '  +2.35%  '

# Row 41: Stellar
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.109'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.05%  '

# Row 42: Maker
$ws.Range("D42").Value = '1.994.33'
$ws.Range("E42").Value = '  +1.57%  '

# Row 43: VeChain
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.0289'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +3.17%  '

# Row 44: ApeXProtocol
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.22'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -7.58%  '

# Row 45: FraxShare
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '10.26'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +4.93%  '

# Row 46: EnergySwap
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '17.76'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.90%  '

# Row 47: NEARProtocol
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.91'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +2.50%  '

# Row 48: MultiversX
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '56.08'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +6.06%  '

# Row 49: RocketPoolETH
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.528.52'
$ws.Range("E49").Value = '  +0.29%  '

# Row 50: Stacks
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.54'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +2.96%  '

# Row 51: THORChain
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '4.51'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.22%  '
